# "Updated status for mathOperator" - Implementation Status.xlsx
#
# - mathOperator row: Implementation Status Partial -> Full, and the Notes
#   text updated from "3 operators unimplemented" to
#   "All operators implemented - most untested".
# - assessmentItem row: Implementation Status "?" -> Full.
# - Re-apply the "Dev issues" AutoFilter on the Implementation Status column
#   (this both drops the now-unused "?" criterion from the filter's value
#   list and recomputes which rows are hidden given the edited values).
# - Move the sheet's selection to E10 (and let Excel drop the stale
#   top-left scroll anchor that pointed at the old selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dev issues")

# mathOperator is row 109: A=Element Name, B=Implementation Status,
# C=Returns, D=..., E=Notes.
$ws.Range("B109").Value = "Full"
$ws.Range("E109").Value = "All operators implemented - most untested"

# assessmentItem is row 10.
$ws.Range("B10").Value = "Full"

# Re-apply the existing AutoFilter criteria on column B (Implementation
# Status, the 2nd column of A:E) now that no row has the value "?" any
# more. This regenerates the filter's value list and re-hides/reveals rows
# to match the edited data (rows 10, 60, 109, 162, 176 and 211 end up
# hidden).
$rng = $ws.Range("A1:E215")
[void]$rng.AutoFilter(2, @("Buggy", "None", "Partial"), 7)

# Move the selection to E10 (clears the stale topLeftCell scroll anchor).
$ws.Activate()
[void]$ws.Range("E10").Select()
